# Apply weekly fruit/vegetable price reshuffle (Achicoria sheet).
# The underlying change is a re-ordering of the daily records (rows 2-14),
# so we simply overwrite the cells that differ between the "before" and
# "after" states, matching the canonical OOXML diff exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 44186

# Row 4
$ws.Range("D4").Value = 44189
$ws.Range("J4").Value = 250

# Row 6
$ws.Range("D6").Value = 44231

# Row 7
$ws.Range("D7").Value = 44188
$ws.Range("J7").Value = 210

# Row 9
$ws.Range("D9").Value = 44208

# Row 10
$ws.Range("D10").Value = 44210
$ws.Range("J10").Value = 340

# Row 11
$ws.Range("D11").Value = 44187
$ws.Range("J11").Value = 160
$ws.Range("K11").Value = 5000
$ws.Range("M11").Value = 5500
$ws.Range("O11").Value = "Provincia de Quillota"
$ws.Range("P11").Value = 344

# Row 12
$ws.Range("D12").Value = 44251
$ws.Range("J12").Value = 120
$ws.Range("L12").Value = 5000
$ws.Range("M12").Value = 5000
$ws.Range("O12").Value = "Región Metropolitana"
$ws.Range("P12").Value = 312

# Row 13
$ws.Range("D13").Value = 44215
$ws.Range("J13").Value = 250

# Row 14
$ws.Range("D14").Value = 44292
$ws.Range("J14").Value = 90
$ws.Range("K14").Value = 6000
$ws.Range("L14").Value = 6000
$ws.Range("M14").Value = 6000
$ws.Range("P14").Value = 375
